$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U4").Value = "A"
$ws.Range("V4").Value = 2.4
$ws.Range("X4").Value = "E=a(KLOC)^b"
$ws.Range("Z4").Value = "Effort"
$ws.Range("AA4").Value = 0.55982240000000005
$ws.Range("AB4").Value = "PM"

$ws.Range("U5").Value = "KLOC"
$ws.Range("V5").Value = 0.25
$ws.Range("X5").Value = "D=c*E^d"
$ws.Range("Z5").Value = "DevTime"
$ws.Range("AA5").Value = 2.0059999999999998
$ws.Range("AB5").Value = "Months"

$ws.Range("U6").Value = "B"
$ws.Range("V6").Value = 1.05
$ws.Range("X6").Value = "SS=E/D"
$ws.Range("Z6").Value = "StaffSize"
$ws.Range("AA6").Value = 0.26900000000000002
$ws.Range("AB6").Value = "Persons"

$ws.Range("U7").Value = "D"
$ws.Range("V7").Value = 0.38
$ws.Range("X7").Value = "P=KLOC/E"
$ws.Range("Z7").Value = "Prod"
$ws.Range("AA7").Value = 0.44600000000000001
$ws.Range("AB7").Value = "KLOC/Person"

$ws.Range("U8").Value = "C"
$ws.Range("V8").Value = 2.5

$ws.Range("U10").Value = ".25^1.05"
$ws.Range("V10").Value = 0.23326
$ws.Range("X10").Value = "E^d = 0.802"

$ws.Range("U11").Value = ".23326 * 2.4"
$ws.Range("V11").Value = 0.55982240000000005

Write-Host "done"
